$d = $word.ActiveDocument

# The engineering-notebook log table currently ends with the 10/2/2022 row.
# Add a new row for 10/5/2022 describing the SRS/SDS progress.
$table = $d.Tables.Item(1)
$newRow = $table.Rows.Add()

$dateCell = $newRow.Cells.Item(1)
$dateCell.Range.Text = "10/5/2022"

$notesCell = $newRow.Cells.Item(2)
$notesCell.Range.Text = "Completed the SRS document for the deliverable of the first sprint.`rXXBLANKLINEXX`rCurrently working on the SDS document for the same deliverable of the first sprint. "

# Turn the placeholder paragraph into a truly empty paragraph (matching the
# blank line already used elsewhere in the notebook) instead of leaving a
# run with empty text behind.
$d.Content.Find.Execute("XXBLANKLINEXX", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
